$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 222.75
$ws.Range("I5").Value = 222.75
$ws.Range("K5").Value = 222.75
$ws.Range("M5").Value = -107.75
$ws.Range("H21").Value = 0
$ws.Range("I21").Value = 0
$ws.Range("K21").Value = 0
$ws.Range("M21").ClearContents()
$ws.Range("H23").Value = 0
$ws.Range("I23").Value = 0
$ws.Range("K23").Value = 0
$ws.Range("M23").ClearContents()
$ws.Range("H62").Value = 4599
$ws.Range("I62").Value = 4497
$ws.Range("K62").Value = 4497
$ws.Range("M62").Value = -3873
$ws.Range("H65").Value = 4599
$ws.Range("I65").Value = 4497
$ws.Range("K65").Value = 22485
$ws.Range("M65").Value = -19365
$ws.Range("H80").Value = 724.375
$ws.Range("J80").Value = 1311.4286
$ws.Range("L80").Value = 3934.2858
$ws.Range("N80").Value = -5930.2858
$ws.Range("H83").Value = 724.375
$ws.Range("J83").Value = 1311.4286
$ws.Range("L83").Value = 11802.8574
$ws.Range("N83").Value = -21786.8574
$ws.Range("H96").Value = 1229.8
$ws.Range("I96").Value = 1514.25
$ws.Range("K96").Value = 4542.75
$ws.Range("M96").Value = -3169.75
$ws.Range("H106").Value = 33999.5
$ws.Range("I106").Value = 33999.5
$ws.Range("K106").Value = 33999.5
$ws.Range("M106").Value = -33368.5

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H11").Value = 15245
$ws.Range("J11").Value = 15245
$ws.Range("L11").Value = 15245
$ws.Range("N11").Value = -15533
$ws.Range("H32").Value = 3057.6365
$ws.Range("I32").Value = 2727.0952
$ws.Range("K32").Value = 2727.0952
$ws.Range("M32").Value = -2440.0952
$ws.Range("H61").Value = 2845.8462
$ws.Range("I61").Value = 2616.25
$ws.Range("K61").Value = 2616.25
$ws.Range("M61").Value = -2404.25
$ws.Range("H74").Value = 2271.6428
$ws.Range("I74").Value = 2254.2307
$ws.Range("K74").Value = 2254.2307
$ws.Range("M74").Value = -1380.2307
$ws.Range("H77").Value = 2271.6428
$ws.Range("I77").Value = 2254.2307
$ws.Range("K77").Value = 11271.1535
$ws.Range("M77").Value = -6903.1535
$ws.Range("H106").Value = 30869.25
$ws.Range("J106").Value = 30869.25
$ws.Range("L106").Value = 30869.25
$ws.Range("N106").Value = -33393.25
$ws.Range("H122").Value = 3720.8
$ws.Range("I122").Value = 3720.8
$ws.Range("K122").Value = 11162.4
$ws.Range("M122").Value = -8712.400000000001
$ws.Range("H125").Value = 91498.75
$ws.Range("J125").Value = 91498.75
$ws.Range("L125").Value = 91498.75
$ws.Range("N125").Value = -101338.75
$ws.Range("H136").Value = 2845.8462
$ws.Range("I136").Value = 2616.25
$ws.Range("K136").Value = 7848.75
$ws.Range("M136").Value = -5298.75

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 333.875
$ws.Range("I22").Value = 340.64285
$ws.Range("J22").Value = 286.5
$ws.Range("K22").Value = 340.64285
$ws.Range("L22").Value = 286.5
$ws.Range("M22").Value = -167.64285
$ws.Range("N22").Value = -632.5
$ws.Range("H86").Value = 2034.7142
$ws.Range("I86").Value = 2226.182
$ws.Range("J86").Value = 1332.6666
$ws.Range("K86").Value = 2226.182
$ws.Range("L86").Value = 1332.6666
$ws.Range("M86").Value = -1103.182
$ws.Range("N86").Value = -3578.6666
$ws.Range("H89").Value = 2034.7142
$ws.Range("I89").Value = 2226.182
$ws.Range("J89").Value = 1332.6666
$ws.Range("K89").Value = 11130.91
$ws.Range("L89").Value = 6663.333000000001
$ws.Range("M89").Value = -5514.91
$ws.Range("N89").Value = -17895.333
$ws.Range("H94").Value = 426.8
$ws.Range("J94").Value = 0
$ws.Range("L94").Value = 0
$ws.Range("N94").ClearContents()
$ws.Range("H107").Value = 0
$ws.Range("I107").Value = 0
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 0
$ws.Range("L107").Value = 0
$ws.Range("M107").ClearContents()
$ws.Range("N107").ClearContents()
$ws.Range("H134").Value = 4479.24
$ws.Range("I134").Value = 4480.3184
$ws.Range("J134").Value = 4471.3335
$ws.Range("K134").Value = 13440.9552
$ws.Range("L134").Value = 13414.0005
$ws.Range("M134").Value = -10905.9552
$ws.Range("N134").Value = -18484.0005

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H18").Value = 40000
$ws.Range("J18").Value = 40000
$ws.Range("L18").Value = 40000
$ws.Range("N18").Value = -40460
$ws.Range("H58").Value = 3351.0557
$ws.Range("I58").Value = 3159.4285
$ws.Range("K58").Value = 3159.4285
$ws.Range("M58").Value = -2956.4285
$ws.Range("H132").Value = 2344.6155
$ws.Range("I132").Value = 2344.6155
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 7033.8465
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -4503.8465
$ws.Range("N132").ClearContents()
$ws.Range("H134").Value = 3195.0833
$ws.Range("I134").Value = 3195.0833
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 9585.249899999999
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = -7050.249899999999
$ws.Range("N134").ClearContents()
$ws.Range("H136").Value = 3351.0557
$ws.Range("I136").Value = 3159.4285
$ws.Range("K136").Value = 9478.2855
$ws.Range("M136").Value = -6928.2855

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 372
$ws.Range("J12").Value = 500
$ws.Range("L12").Value = 1500
$ws.Range("N12").Value = -1846
$ws.Range("H34").Value = 45829.5
$ws.Range("I34").Value = 685
$ws.Range("J34").Value = 52278.715
$ws.Range("K34").Value = 2055
$ws.Range("L34").Value = 156836.145
$ws.Range("M34").Value = -1971
$ws.Range("N34").Value = -157004.145
$ws.Range("H131").Value = 1569.6666
$ws.Range("J131").Value = 0
$ws.Range("L131").Value = 0
$ws.Range("N131").ClearContents()

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 642.5714
$ws.Range("I97").Value = 616.3333
$ws.Range("K97").Value = 616.3333
$ws.Range("M97").Value = -120.3333
$ws.Range("H122").Value = 1865.1
$ws.Range("I122").Value = 1949.2222
$ws.Range("J122").Value = 1108
$ws.Range("K122").Value = 5847.6666
$ws.Range("L122").Value = 3324
$ws.Range("M122").Value = -3397.6666
$ws.Range("N122").Value = -8224
$ws.Range("H126").Value = 2562.375
$ws.Range("I126").Value = 2642.8572
$ws.Range("K126").Value = 7928.571599999999
$ws.Range("M126").Value = -5458.571599999999

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 9124.625
$ws.Range("I7").Value = 9213.143
$ws.Range("K7").Value = 9213.143
$ws.Range("M7").Value = -9101.143
$ws.Range("H63").Value = 51750
$ws.Range("J63").Value = 51750
$ws.Range("L63").Value = 51750
$ws.Range("N63").Value = -53248
$ws.Range("H66").Value = 51750
$ws.Range("J66").Value = 51750
$ws.Range("L66").Value = 155250
$ws.Range("N66").Value = -162738
$ws.Range("H126").Value = 9124.625
$ws.Range("I126").Value = 9213.143
$ws.Range("K126").Value = 27639.429
$ws.Range("M126").Value = -25169.429
$ws.Range("H132").Value = 4909.9473
$ws.Range("I132").Value = 4393.125
$ws.Range("K132").Value = 13179.375
$ws.Range("M132").Value = -10649.375
$ws.Range("H136").Value = 2290.3333
$ws.Range("I136").Value = 1866.2858
$ws.Range("K136").Value = 5598.857400000001
$ws.Range("M136").Value = -3048.857400000001

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 1051.5
$ws.Range("I126").Value = 1051.5
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 3154.5
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -684.5
$ws.Range("N126").ClearContents()
$ws.Range("H136").Value = 5315.343
$ws.Range("I136").Value = 1724.3334
$ws.Range("K136").Value = 5173.0002
$ws.Range("M136").Value = -2624.0002
